# 案件情報.xlsx - "ランサーズ" sheet refresh: append a new 2025-10-21 18:24:52 scrape run.
# Every visible row is rewritten in place with the new snapshots values; the two
# newest postings ("出会い系アプリ" and "SESエンジニア") push the list from 10 rows to 12,
# so the sheet grows from A1:H11 to A1:H13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlinks before rewriting F2:F13 so they can be rebuilt
# cleanly against the new row order (Hyperlinks.Delete on any cell in the sheet
# clears the whole collection in this engine).
$ws.Range("F2").Hyperlinks.Delete()

# Row 2: 【急募】ebayAPIを活用したShippingポリシー設定の専門家募集
$ws.Range("A2").Value = '2025-10-21 18:24:52'
$ws.Range("B2").Value = '【急募】ebayAPIを活用したShippingポリシー設定の専門家募集'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5415908'
$ws.Range("G2").Value = 183
$ws.Range("H2").Value = '🔥API'

# Row 3: 【 急募! 】 JS、PHPを使用したWEBシステムの開発、修正の対応
$ws.Range("A3").Value = '2025-10-21 18:24:52'
$ws.Range("B3").Value = '【 急募! 】 JS、PHPを使用したWEBシステムの開発、修正の対応'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5417295'
$ws.Range("G3").Value = 110
$ws.Range("H3").Value = '◆開発 ○PHP'

# Row 4: システムの開発補助や運営サポート【フルリモート×長期】
$ws.Range("A4").Value = '2025-10-21 18:24:52'
$ws.Range("B4").Value = 'システムの開発補助や運営サポート【フルリモート×長期】'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5408664'
$ws.Range("G4").Value = 90
$ws.Range("H4").Value = '◆開発'

# Row 5: MySQLバージョンアップ(ロリポップ/WordPress/1データベースに8サイト)
$ws.Range("A5").Value = '2025-10-21 18:24:52'
$ws.Range("B5").Value = 'MySQLバージョンアップ(ロリポップ/WordPress/1データベースに8サイト)'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5417433'
$ws.Range("G5").Value = 88
$ws.Range("H5").Value = '◇MySQL ○WordPress'

# Row 6: 【急募】MT4/MT5用FX自動売買システムの開発者募集
$ws.Range("A6").Value = '2025-10-21 18:24:52'
$ws.Range("B6").Value = '【急募】MT4/MT5用FX自動売買システムの開発者募集'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5417377'
$ws.Range("G6").Value = 83
$ws.Range("H6").Value = '◆開発'

# Row 7: 【急募】Salesforce・MA・CRMコンサルタント経験者を探しています!
$ws.Range("A7").Value = '2025-10-21 18:24:52'
$ws.Range("B7").Value = '【急募】Salesforce・MA・CRMコンサルタント経験者を探しています!'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5371747'
$ws.Range("G7").Value = 48
$ws.Range("H7").Value = '◆コンサル'

# Row 8: 【急募】出会い系アプリの制作から運用までお任せ!
$ws.Range("A8").Value = '2025-10-21 18:24:52'
$ws.Range("B8").Value = '【急募】出会い系アプリの制作から運用までお任せ!'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5417817'
$ws.Range("G8").Value = 45
$ws.Range("H8").Value = '◇アプリ'

# Row 9: 【高齢者支援】見守りアプリとマニュアル制作の依頼
$ws.Range("A9").Value = '2025-10-21 18:24:52'
$ws.Range("B9").Value = '【高齢者支援】見守りアプリとマニュアル制作の依頼'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5417267'
$ws.Range("G9").Value = 45
$ws.Range("H9").Value = '◇アプリ'

# Row 10: 【 急募】コミュニティサイトのカスタマイズ、修正
$ws.Range("A10").Value = '2025-10-21 18:24:52'
$ws.Range("B10").Value = '【 急募】コミュニティサイトのカスタマイズ、修正'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5417308'
$ws.Range("G10").Value = 38
$ws.Range("H10").Value = '◇サイト'

# Row 11: 【Webarena suiteX/DNS】ドメイン設定変更によるウェブサイト分割とサイト切り替え
$ws.Range("A11").Value = '2025-10-21 18:24:52'
$ws.Range("B11").Value = '【Webarena suiteX/DNS】ドメイン設定変更によるウェブサイト分割とサイト切り替え'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5417544'
$ws.Range("G11").Value = 30
$ws.Range("H11").Value = '◇サイト'

# Row 12: 【フルリモート】SESエンジニア募集|スキルに応じて30〜40万円/月|複数案件あり・継続前提
$ws.Range("A12").Value = '2025-10-21 18:24:52'
$ws.Range("B12").Value = '【フルリモート】SESエンジニア募集|スキルに応じて30〜40万円/月|複数案件あり・継続前提'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5417644'
$ws.Range("G12").Value = 25

# Row 13: 【データ加工のプロ募集】施設情報データの修正・整備依頼
$ws.Range("A13").Value = '2025-10-21 18:24:52'
$ws.Range("B13").Value = '【データ加工のプロ募集】施設情報データの修正・整備依頼'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5417622'
$ws.Range("G13").Value = 10

# Re-create the hyperlinks for the URL column, in row order, keeping the same
# "Hyperlink" cell style (s="1") that the sheet already used for column F.
$urls = @(
  'https://www.lancers.jp/work/detail/5415908',
  'https://www.lancers.jp/work/detail/5417295',
  'https://www.lancers.jp/work/detail/5408664',
  'https://www.lancers.jp/work/detail/5417433',
  'https://www.lancers.jp/work/detail/5417377',
  'https://www.lancers.jp/work/detail/5371747',
  'https://www.lancers.jp/work/detail/5417817',
  'https://www.lancers.jp/work/detail/5417267',
  'https://www.lancers.jp/work/detail/5417308',
  'https://www.lancers.jp/work/detail/5417544',
  'https://www.lancers.jp/work/detail/5417644',
  'https://www.lancers.jp/work/detail/5417622'
)
for ($i = 0; $i -lt $urls.Length; $i++) {
  $row = $i + 2
  $cell = $ws.Cells.Item($row, 6)
  $ws.Hyperlinks.Add($cell, $urls[$i])
  $cell.Style = "Hyperlink"
}

Write-Output "applied 2025-10-21 18:24 refresh: A1:H13"
